# Add "N. crash:" (count of True values), "Mean:" and "Standard Deviation:"
# summary rows below the data table (rows 53-55), per columns G..P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 - N. crash: count of "True" crash values in column G
$ws.Range("F53").Value = "N. crash:"
$ws.Range("G53").Formula = '=COUNTIF(G2:G51,"True")'

# Row 54 - Mean: average of each numeric column H..P
$ws.Range("F54").Value = "Mean:"
$ws.Range("H54").Formula = "=AVERAGE(H2:H51)"
$ws.Range("I54").Formula = "=AVERAGE(I2:I51)"
$ws.Range("J54").Formula = "=AVERAGE(J2:J51)"
$ws.Range("K54").Formula = "=AVERAGE(K2:K51)"
$ws.Range("L54").Formula = "=AVERAGE(L2:L51)"
$ws.Range("M54").Formula = "=AVERAGE(M2:M51)"
$ws.Range("N54").Formula = "=AVERAGE(N2:N51)"
$ws.Range("O54").Formula = "=AVERAGE(O2:O51)"
$ws.Range("P54").Formula = "=AVERAGE(P2:P51)"

# Row 55 - Standard Deviation: sample standard deviation of each numeric column H..P
$ws.Range("F55").Value = "Standard Deviation:"
$ws.Range("H55").Formula = "=STDEV.S(H2:H51)"
$ws.Range("I55").Formula = "=STDEV.S(I2:I51)"
$ws.Range("J55").Formula = "=STDEV.S(J2:J51)"
$ws.Range("K55").Formula = "=STDEV.S(K2:K51)"
$ws.Range("L55").Formula = "=STDEV.S(L2:L51)"
$ws.Range("M55").Formula = "=STDEV.S(M2:M51)"
$ws.Range("N55").Formula = "=STDEV.S(N2:N51)"
$ws.Range("O55").Formula = "=STDEV.S(O2:O51)"
$ws.Range("P55").Formula = "=STDEV.S(P2:P51)"

# Keep the current selection where the author left it when they last saved
$null = $ws.Range("H61").Select()
